$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.538.19'
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("D3").Value = '1.912.67'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.633'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.74%  '
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.14'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0707'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0997'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").Value = '2.187.04'
$ws.Range("E12").Value = '  +3.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '12.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.49%  '
$ws.Range("D14").Value = '1.904.95'
$ws.Range("E14").Value = '  +2.95%  '
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("E16").Value = '  +3.94%  '
$ws.Range("D17").Value = '35.572.31'
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '72.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.08%  '
$ws.Range("D19").Value = '0.0₃0821'
$ws.Range("E19").Value = '  +3.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '243.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.76%  '
$ws.Range("E22").Value = '  +2.19%  '
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("E26").Value = '  +20.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.34%  '
$ws.Range("E29").Value = '  +0.87%  '
$ws.Range("E30").Value = '  +26.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0572'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.90%  '
$ws.Range("E32").Value = '  +3.50%  '
$ws.Range("E33").Value = '  +5.85%  '
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("E35").Value = '  +6.32%  '
$ws.Range("E36").Value = '  +12.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("E38").Value = '  +3.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0206'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '91.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.36%  '
$ws.Range("D42").Value = '1.355.47'
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +43.01%  '
$ws.Range("E44").Value = '  +13.32%  '
$ws.Range("E45").Value = '  +2.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("E48").Value = '  +0.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.65%  '
$ws.Range("D50").Value = '2.096.86'
$ws.Range("E50").Value = '  +3.36%  '
$ws.Range("E51").Value = '  +2.07%  '
